$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 23.699655
$ws.Range("H2").Value = 71.09896499999999
$ws.Range("I2").Value = 0.4841969272415696
$ws.Range("J2").Value = 0.4841969272415697
$ws.Range("M2").Value = 2.724001666666667
$ws.Range("N2").Value = 8.172005
$ws.Range("O2").Value = 0.04635500474236593
$ws.Range("P2").Value = 0.04635500474236593
$ws.Range("Q2").Value = 64.557899719425
$ws.Range("R2").Value = 581.021097474825
$ws.Range("S2").Value = 0.02244495085852197
$ws.Range("T2").Value = 0.02244495085852197
$ws.Range("G3").Value = 23.699655
$ws.Range("H3").Value = 71.09896499999999
$ws.Range("I3").Value = 0.4841969272415696
$ws.Range("J3").Value = 0.4841969272415697
$ws.Range("O3").Value = 0.6912512390256352
$ws.Range("P3").Value = 0.6912512390256351
$ws.Range("Q3").Value = 962.6949326824199
$ws.Range("R3").Value = 8664.254394141781
$ws.Range("S3").Value = 0.3347017258881403
$ws.Range("T3").Value = 0.3347017258881403
$ws.Range("G4").Value = 23.699655
$ws.Range("H4").Value = 71.09896499999999
$ws.Range("I4").Value = 0.4841969272415696
$ws.Range("J4").Value = 0.4841969272415697
$ws.Range("M4").Value = 15.419285
$ws.Range("N4").Value = 46.257855
$ws.Range("O4").Value = 0.2623937562319988
$ws.Range("P4").Value = 0.2623937562319988
$ws.Range("Q4").Value = 365.4317348466749
$ws.Range("R4").Value = 3288.885613620075
$ws.Range("S4").Value = 0.1270502504949073
$ws.Range("T4").Value = 0.1270502504949073
$ws.Range("H5").Value = 58.032849
$ws.Range("I5").Value = 0.3952142927098025
$ws.Range("J5").Value = 0.3952142927098025
$ws.Range("M5").Value = 2.724001666666667
$ws.Range("N5").Value = 8.172005
$ws.Range("O5").Value = 0.04635500474236593
$ws.Range("P5").Value = 0.04635500474236593
$ws.Range("Q5").Value = 52.69385913247167
$ws.Range("R5").Value = 474.244732192245
$ws.Range("S5").Value = 0.01832016041281369
$ws.Range("T5").Value = 0.01832016041281369
$ws.Range("H6").Value = 58.032849
$ws.Range("I6").Value = 0.3952142927098025
$ws.Range("J6").Value = 0.3952142927098025
$ws.Range("O6").Value = 0.6912512390256352
$ws.Range("P6").Value = 0.6912512390256351
$ws.Range("Q6").Value = 785.7769752544788
$ws.Range("R6").Value = 7071.992777290308
$ws.Range("S6").Value = 0.2731923695162911
$ws.Range("T6").Value = 0.273192369516291
$ws.Range("H7").Value = 58.032849
$ws.Range("I7").Value = 0.3952142927098025
$ws.Range("J7").Value = 0.3952142927098025
$ws.Range("M7").Value = 15.419285
$ws.Range("N7").Value = 46.257855
$ws.Range("O7").Value = 0.2623937562319988
$ws.Range("P7").Value = 0.2623937562319988
$ws.Range("Q7").Value = 298.275012697655
$ws.Range("R7").Value = 2684.475114278895
$ws.Range("S7").Value = 0.1037017627806977
$ws.Range("T7").Value = 0.1037017627806977
$ws.Range("G8").Value = 5.902376333333333
$ws.Range("H8").Value = 17.707129
$ws.Range("I8").Value = 0.1205887800486278
$ws.Range("J8").Value = 0.1205887800486278
$ws.Range("M8").Value = 2.724001666666667
$ws.Range("N8").Value = 8.172005
$ws.Range("O8").Value = 0.04635500474236593
$ws.Range("P8").Value = 0.04635500474236593
$ws.Range("Q8").Value = 16.07808296929389
$ws.Range("R8").Value = 144.702746723645
$ws.Range("S8").Value = 0.005589893471030265
$ws.Range("T8").Value = 0.005589893471030265
$ws.Range("G9").Value = 5.902376333333333
$ws.Range("H9").Value = 17.707129
$ws.Range("I9").Value = 0.1205887800486278
$ws.Range("J9").Value = 0.1205887800486278
$ws.Range("O9").Value = 0.6912512390256352
$ws.Range("P9").Value = 0.6912512390256351
$ws.Range("Q9").Value = 239.7582490920076
$ws.Range("R9").Value = 2157.824241828068
$ws.Range("S9").Value = 0.08335714362120379
$ws.Range("T9").Value = 0.08335714362120378
$ws.Range("G10").Value = 5.902376333333333
$ws.Range("H10").Value = 17.707129
$ws.Range("I10").Value = 0.1205887800486278
$ws.Range("J10").Value = 0.1205887800486278
$ws.Range("M10").Value = 15.419285
$ws.Range("N10").Value = 46.257855
$ws.Range("O10").Value = 0.2623937562319988
$ws.Range("P10").Value = 0.2623937562319988
$ws.Range("Q10").Value = 91.01042286092166
$ws.Range("R10").Value = 819.0938057482949
$ws.Range("S10").Value = 0.03164174295639378
$ws.Range("T10").Value = 0.03164174295639378
